$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.115.65'

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -4.57%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.653.17'

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -3.48%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.82'

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -3.96%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5096'

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -3.71%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.12%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.21%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06413'

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -4.34%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.91'

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.75%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07791'

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.48%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.652.06'

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.56%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.280'

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -5.19%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.881.51'

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.42%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5520'

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -5.42%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8011'

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.66%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.94'

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -6.25%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.149.54'

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -4.45%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.20%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '209.32'

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -6.29%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.416'

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -4.66%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.15%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.040'

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.39%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.12%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.54'

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.00%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +3.26%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.41%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.973'

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.73%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.82'

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.12%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05091'

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -5.09%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.242'

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.90%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.340'

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.215'

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -6.47%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.566'

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.34%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.750'

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.30%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.366'

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.31%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9261'

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.66%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.167.57'

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +6.57%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5677'

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.09%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01591'

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.76%  '

$ws.Range("B41").Value = 'mCoin'

$ws.Range("C41").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.556'

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.58%  '

$ws.Range("B42").Value = 'PaxDollar'

$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.006'

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.11%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8319'

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.26%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.665'

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.22%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '100.39'

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.69%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.791.60'

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.42%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₈117'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.80%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4555'

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.39%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '55.66'

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.66%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.007'

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.33%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.864'

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.09%  '
